$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B10").Value = "pass"
$ws.Range("A10").Value = "completing a horizontal set of 4. The AI sets a 2 on the right of this group."
$ws.Range("A12").Value = "AIVerticalTest1(): This test sets up the board in such a manner that the AI opponent places a 2 on any stack of 3 1's that the human player places."
$ws.Range("A9").Value = "AIHorizontalTest1(): This test sets up the board in a manner such that the AI opponent places a piece to prevent the human player from winning by "
$ws.Range("A14").Value = "AIDiagonalRightTest1(): This test sets up the board in such a manner that the AI opponent places a 2 on a group of 1's that form a ""right diagonal"". "
$ws.Range("A16").Value = "AIDiagonalLeftTest1(): This test sets up the board in such a manner that the AI opponent places a 2 on a group of 1's that form a ""left diagonal""."

$ws.Range("B12").Value = "pass"
$ws.Range("B14").Value = "pass"
$ws.Range("B16").Value = "pass"

$ws.Range("B16").Select()
